# "Add and fix skills" - append two new skill rows to the table, grow the
# table range/autofilter to match, and leave the selection on the newly
# added data (A23, matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new skills below the existing data (rows 1-28 -> 1-30).
$ws.Range("A29").Value = "BASIC"
$ws.Range("A30").Value = "Email"

# Grow Table1 so the new rows are included (ref + autoFilter both extend
# from A1:A28 to A1:A30).
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A1:A30"))

# Match the saved selection from the diff.
[void]$ws.Range("A23").Select()
